# Added periodic & upfront related scenarios
#
# The "repaymentstrategy" value on the ProductLoanInput sheet (cell B17)
# is changed from "RBI (India)" to "Overdue/Due Fee/Int,Principal", and
# the ProductLoanInput sheet/cell B17 becomes the active selection
# (instead of ProductLoanOutput!B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

$ws.Activate()
$ws.Range("B17").Select()
